$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.284.74'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.866.49'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.67'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4701'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06564'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.36'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07826'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.68'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.872.79'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6965'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.088'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.02'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.194.82'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.80'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007651'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.53%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.105.28'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.237'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.169'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.459'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.55'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.938'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.365'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09905'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.355'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.458'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.047'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04722'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7025'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01874'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.756'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.310'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.95'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.949'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4170'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8364'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.01'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '972.62'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.098'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.112'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05677'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.27%  '
